$d = $word.ActiveDocument

$d.Content.Find.Execute("47×54=2538", $true, $false, $false, $false, $false, $true, 1, $false, "76×43=3268", 2) | Out-Null
$d.Content.Find.Execute("66×89=5874", $true, $false, $false, $false, $false, $true, 1, $false, "26×93=2418", 2) | Out-Null
$d.Content.Find.Execute("37×16=592", $true, $false, $false, $false, $false, $true, 1, $false, "37×52=1924", 2) | Out-Null
$d.Content.Find.Execute("94×87=8178", $true, $false, $false, $false, $false, $true, 1, $false, "59×47=2773", 2) | Out-Null
$d.Content.Find.Execute("48×25=1200", $true, $false, $false, $false, $false, $true, 1, $false, "21×98=2058", 2) | Out-Null
$d.Content.Find.Execute("23×61=1403", $true, $false, $false, $false, $false, $true, 1, $false, "49×55=2695", 2) | Out-Null
$d.Content.Find.Execute("46×42=1932", $true, $false, $false, $false, $false, $true, 1, $false, "90×44=3960", 2) | Out-Null
$d.Content.Find.Execute("36×38=1368", $true, $false, $false, $false, $false, $true, 1, $false, "39×71=2769", 2) | Out-Null
$d.Content.Find.Execute("27×11=297", $true, $false, $false, $false, $false, $true, 1, $false, "50×46=2300", 2) | Out-Null
$d.Content.Find.Execute("82×50=4100", $true, $false, $false, $false, $false, $true, 1, $false, "97×24=2328", 2) | Out-Null
$d.Content.Find.Execute("81×63=5103", $true, $false, $false, $false, $false, $true, 1, $false, "71×64=4544", 2) | Out-Null
$d.Content.Find.Execute("74×27=1998", $true, $false, $false, $false, $false, $true, 1, $false, "61×91=5551", 2) | Out-Null
$d.Content.Find.Execute("55×46=2530", $true, $false, $false, $false, $false, $true, 1, $false, "92×40=3680", 2) | Out-Null
$d.Content.Find.Execute("79×27=2133", $true, $false, $false, $false, $false, $true, 1, $false, "62×75=4650", 2) | Out-Null
$d.Content.Find.Execute("29×98=2842", $true, $false, $false, $false, $false, $true, 1, $false, "87×77=6699", 2) | Out-Null
$d.Content.Find.Execute("11×53=583", $true, $false, $false, $false, $false, $true, 1, $false, "86×16=1376", 2) | Out-Null
$d.Content.Find.Execute("82×78=6396", $true, $false, $false, $false, $false, $true, 1, $false, "77×36=2772", 2) | Out-Null
$d.Content.Find.Execute("71×76=5396", $true, $false, $false, $false, $false, $true, 1, $false, "13×87=1131", 2) | Out-Null
$d.Content.Find.Execute("38×17=646", $true, $false, $false, $false, $false, $true, 1, $false, "87×14=1218", 2) | Out-Null
$d.Content.Find.Execute("72×47=3384", $true, $false, $false, $false, $false, $true, 1, $false, "79×59=4661", 2) | Out-Null
$d.Content.Find.Execute("85×28=2380", $true, $false, $false, $false, $false, $true, 1, $false, "46×62=2852", 2) | Out-Null
$d.Content.Find.Execute("77×90=6930", $true, $false, $false, $false, $false, $true, 1, $false, "29×43=1247", 2) | Out-Null
$d.Content.Find.Execute("11×67=737", $true, $false, $false, $false, $false, $true, 1, $false, "16×75=1200", 2) | Out-Null
$d.Content.Find.Execute("31×23=713", $true, $false, $false, $false, $false, $true, 1, $false, "24×22=528", 2) | Out-Null
$d.Content.Find.Execute("49×52=2548", $true, $false, $false, $false, $false, $true, 1, $false, "19×71=1349", 2) | Out-Null
